# Update "want to go" counts (column F) on several rows across sheets,
# mirroring a refreshed scrape of event attendance numbers.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F4").Value = 7737
$wsExpo.Range("F5").Value = 5632
$wsExpo.Range("F10").Value = 253
$wsExpo.Range("F11").Value = 233

# Sheet "演出" (Performance)
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 85

# Sheet "全部类型" (All Types) - aggregated view of all events
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 7737
$wsAll.Range("F5").Value = 5632
$wsAll.Range("F10").Value = 253
$wsAll.Range("F11").Value = 85
$wsAll.Range("F13").Value = 233
